$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is currently blank (the sheet data jumps from row 1 straight to row 3).
# Deleting that blank row shifts every data row up by one, which is the net
# effect captured by the diff (old row 3 "das/Kind/Child/Family" becomes the
# new row 2, ... old row 30 "die/Sommersprossen/Freckles/Face" becomes the
# new (and now last) row 29).
$ws.Rows.Item(2).Delete()

# Restore the selection to match the edited workbook (row 2 selected).
$ws.Range("A2:XFD2").Select()
